$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.065017227438009
$ws.Range("D2").Value = 1.064521878082556
$ws.Range("E2").Value = 1.069321204400903
$ws.Range("F2").Value = 1.077885431701603
$ws.Range("I2").Value = 1.045280637023585
$ws.Range("J2").Value = 1.069974499192607
$ws.Range("K2").Value = 1.067237589073905
$ws.Range("L2").Value = 1.072024024188067
$ws.Range("M2").Value = 1.080565555766972
$ws.Range("N2").Value = 1.071493985709251
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.066456007069047
$ws.Range("D3").Value = 1.065600904426902
$ws.Range("E3").Value = 1.070561642044949
$ws.Range("F3").Value = 1.079113077715598
$ws.Range("I3").Value = 1.045590827535343
$ws.Range("J3").Value = 1.071066400562702
$ws.Range("K3").Value = 1.068131176136955
$ws.Range("L3").Value = 1.073079557480332
$ws.Range("M3").Value = 1.081609979979161
$ws.Range("N3").Value = 1.072587437704535
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.067386365947234
$ws.Range("D4").Value = 1.066298355658953
$ws.Range("E4").Value = 1.071363927145318
$ws.Range("F4").Value = 1.079907018369774
$ws.Range("I4").Value = 1.045789877815952
$ws.Range("J4").Value = 1.071771845065586
$ws.Range("K4").Value = 1.068708038694126
$ws.Range("L4").Value = 1.073761631195808
$ws.Range("M4").Value = 1.082284792841264
$ws.Range("N4").Value = 1.073293884019528
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.067777343295834
$ws.Range("D5").Value = 1.066591386992807
$ws.Range("E5").Value = 1.071701124580298
$ws.Range("F5").Value = 1.080240691467696
$ws.Range("I5").Value = 1.045873161244281
$ws.Range("J5").Value = 1.07206815654552
$ws.Range("K5").Value = 1.068950231503036
$ws.Range("L5").Value = 1.074048156036818
$ws.Range("M5").Value = 1.0825682476869
$ws.Range("N5").Value = 1.073590616295751
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.06784298171759
$ws.Range("D6").Value = 1.066640577930054
$ws.Range("E6").Value = 1.071757736726115
$ws.Range("F6").Value = 1.080296710898442
$ws.Range("I6").Value = 1.045887121628585
$ws.Range("J6").Value = 1.072117893547569
$ws.Range("K6").Value = 1.068990878036818
$ws.Range("L6").Value = 1.07409625206477
$ws.Range("M6").Value = 1.082615827216131
$ws.Range("N6").Value = 1.073640423930047
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.067391590771593
$ws.Range("D7").Value = 1.066302271850868
$ws.Range("E7").Value = 1.071368433116974
$ws.Range("F7").Value = 1.079911477313483
$ws.Range("I7").Value = 1.04579099221313
$ws.Range("J7").Value = 1.071775805401421
$ws.Range("K7").Value = 1.068711276141393
$ws.Range("L7").Value = 1.073765460611415
$ws.Range("M7").Value = 1.082288581304544
$ws.Range("N7").Value = 1.073297849979494
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.065503602600475
$ws.Range("D8").Value = 1.064886697159704
$ws.Range("E8").Value = 1.069740492789738
$ws.Range("F8").Value = 1.078300410938949
$ws.Range("I8").Value = 1.045385812440761
$ws.Range("J8").Value = 1.070343739752837
$ws.Range("K8").Value = 1.067539861637087
$ws.Range("L8").Value = 1.072380940309811
$ws.Range("M8").Value = 1.080918732441391
$ws.Range("N8").Value = 1.071863750633429
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.062171704826846
$ws.Range("D9").Value = 1.062386404283303
$ws.Range("E9").Value = 1.066868945097505
$ws.Range("F9").Value = 1.075458089717916
$ws.Range("I9").Value = 1.044659047186236
$ws.Range("J9").Value = 1.067811773073609
$ws.Range("K9").Value = 1.065465242407377
$ws.Range("L9").Value = 1.069934006832449
$ws.Range("M9").Value = 1.07849710382121
$ws.Range("N9").Value = 1.069328188271096
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.059946763193781
$ws.Range("D10").Value = 1.060715428827481
$ws.Range("E10").Value = 1.064952420622922
$ws.Range("F10").Value = 1.073560720224644
$ws.Range("I10").Value = 1.044165876392784
$ws.Range("J10").Value = 1.066117885058849
$ws.Range("K10").Value = 1.064074985559167
$ws.Range("L10").Value = 1.068297667825993
$ws.Range("M10").Value = 1.076877276469491
$ws.Range("N10").Value = 1.067631894741063
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.058982399529594
$ws.Range("D11").Value = 1.059990864360917
$ws.Range("E11").Value = 1.064121986196017
$ws.Range("F11").Value = 1.07273850460689
$ws.Range("I11").Value = 1.043950259285425
$ws.Range("J11").Value = 1.065382964545046
$ws.Range("K11").Value = 1.063471250381336
$ws.Range("L11").Value = 1.067587874755208
$ws.Range("M11").Value = 1.07617455038904
$ws.Range("N11").Value = 1.066895930555834
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.05862404342732
$ws.Range("D12").Value = 1.059721572086583
$ws.Range("E12").Value = 1.06381343668542
$ws.Range("F12").Value = 1.072432996815136
$ws.Range("I12").Value = 1.043869857048847
$ws.Range("J12").Value = 1.065109759529659
$ws.Range("K12").Value = 1.063246730907779
$ws.Range("L12").Value = 1.067324034843306
$ws.Range("M12").Value = 1.075913323239736
$ws.Range("N12").Value = 1.066622337557988
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.058700918842414
$ws.Range("D13").Value = 1.059779343342025
$ws.Range("E13").Value = 1.063879625704152
$ws.Range("F13").Value = 1.072498533874211
$ws.Range("I13").Value = 1.04388711776138
$ws.Range("J13").Value = 1.065168373081646
$ws.Range("K13").Value = 1.06329490315631
$ws.Range("L13").Value = 1.067380638110817
$ws.Range("M13").Value = 1.075969366611327
$ws.Range("N13").Value = 1.066681034347942
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.058952780765867
$ws.Range("D14").Value = 1.059968607795129
$ws.Range("E14").Value = 1.064096483240175
$ws.Range("F14").Value = 1.072713253305907
$ws.Range("I14").Value = 1.043943619592826
$ws.Range("J14").Value = 1.065360385902346
$ws.Range("K14").Value = 1.063452696963267
$ws.Range("L14").Value = 1.067566069587149
$ws.Range("M14").Value = 1.076152961441203
$ws.Range("N14").Value = 1.066873319848872
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.05910794135531
$ws.Range("D15").Value = 1.060085198996559
$ws.Range("E15").Value = 1.064230084385265
$ws.Range("F15").Value = 1.07284553562835
$ws.Range("I15").Value = 1.043978390794884
$ws.Range("J15").Value = 1.065478661690536
$ws.Range("K15").Value = 1.063549883678081
$ws.Range("L15").Value = 1.067680294596634
$ws.Range("M15").Value = 1.0762660532169
$ws.Range("N15").Value = 1.066991763602247
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.06001074433281
$ws.Range("D16").Value = 1.060763493985261
$ws.Range("E16").Value = 1.065007521536212
$ws.Range("F16").Value = 1.073615274112207
$ws.Range("I16").Value = 1.044180142446692
$ws.Range("J16").Value = 1.066166628291537
$ws.Range("K16").Value = 1.064115016440113
$ws.Range("L16").Value = 1.068344747852533
$ws.Range("M16").Value = 1.076923885757473
$ws.Range("N16").Value = 1.067680707194732
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.060576790500691
$ws.Range("D17").Value = 1.061188694772217
$ws.Range("E17").Value = 1.065495032750744
$ws.Range("F17").Value = 1.07409793611743
$ws.Range("I17").Value = 1.044306140564179
$ws.Range("J17").Value = 1.066597778915736
$ws.Range("K17").Value = 1.064469039899362
$ws.Range("L17").Value = 1.068761205203377
$ws.Range("M17").Value = 1.077336168108797
$ws.Range("N17").Value = 1.068112470102272
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.06090686446417
$ws.Range("D18").Value = 1.061436608908673
$ws.Range("E18").Value = 1.06577933539337
$ws.Range("F18").Value = 1.074379403006219
$ws.Range("I18").Value = 1.044379433412178
$ws.Range("J18").Value = 1.066849121592768
$ws.Range("K18").Value = 1.064675367665563
$ws.Range("L18").Value = 1.069003997620911
$ws.Range("M18").Value = 1.077576517331142
$ws.Range("N18").Value = 1.068364169714734
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.061019395865715
$ws.Range("D19").Value = 1.061521124652177
$ws.Range("E19").Value = 1.065876266141974
$ws.Range("F19").Value = 1.074475365592551
$ws.Range("I19").Value = 1.044404390542375
$ws.Range("J19").Value = 1.066934799295234
$ws.Range("K19").Value = 1.064745691680989
$ws.Range("L19").Value = 1.069086763300865
$ws.Range("M19").Value = 1.077658448583993
$ws.Range("N19").Value = 1.068449969089365
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.060516068580253
$ws.Range("D20").Value = 1.061143084943833
$ws.Range("E20").Value = 1.065442733068745
$ws.Range("F20").Value = 1.074046157462349
$ws.Range("I20").Value = 1.044292642826805
$ws.Range("J20").Value = 1.066551535065273
$ws.Range("K20").Value = 1.06443107394721
$ws.Range("L20").Value = 1.068716535722688
$ws.Range("M20").Value = 1.077291947390972
$ws.Range("N20").Value = 1.068066160580236
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.058878617853607
$ws.Range("D21").Value = 1.059912878485749
$ws.Range("E21").Value = 1.064032626595566
$ws.Range("F21").Value = 1.07265002658474
$ws.Range("I21").Value = 1.043926989847267
$ws.Range("J21").Value = 1.065303849096189
$ws.Range("K21").Value = 1.063406237970801
$ws.Range("L21").Value = 1.067511469946285
$ws.Range("M21").Value = 1.076098902983109
$ws.Range("N21").Value = 1.066816702753965
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.05784822548689
$ws.Range("D22").Value = 1.059138490152263
$ws.Range("E22").Value = 1.063145518643713
$ws.Range("F22").Value = 1.071771641457792
$ws.Range("I22").Value = 1.043695280954021
$ws.Range("J22").Value = 1.0645180882296
$ws.Range("K22").Value = 1.062760346095205
$ws.Range("L22").Value = 1.066752689391985
$ws.Range("M22").Value = 1.075347610061107
$ws.Range("N22").Value = 1.066029826016818
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.05839453966496
$ws.Range("D23").Value = 1.059549095323443
$ws.Range("E23").Value = 1.063615841816378
$ws.Range("F23").Value = 1.072237346506941
$ws.Range("I23").Value = 1.043818286073557
$ws.Range("J23").Value = 1.064934758692391
$ws.Range("K23").Value = 1.063102892348974
$ws.Range("L23").Value = 1.067155039587934
$ws.Range("M23").Value = 1.075745997558495
$ws.Range("N23").Value = 1.066447088199458
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.06054350648376
$ws.Range("D24").Value = 1.061163694365592
$ws.Range("E24").Value = 1.065466365211736
$ws.Range("F24").Value = 1.074069554198259
$ws.Range("I24").Value = 1.044298742490378
$ws.Range("J24").Value = 1.066572431103697
$ws.Range("K24").Value = 1.0644482296461
$ws.Range("L24").Value = 1.068736720309172
$ws.Range("M24").Value = 1.07731192922554
$ws.Range("N24").Value = 1.068087086293432
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.063033706169958
$ws.Range("D25").Value = 1.06303350263614
$ws.Range("E25").Value = 1.06761167461338
$ws.Range("F25").Value = 1.07619332241316
$ws.Range("I25").Value = 1.044848455236036
$ws.Range("J25").Value = 1.068467373700393
$ws.Range("K25").Value = 1.066002834481693
$ws.Range("L25").Value = 1.070567474064155
$ws.Range("M25").Value = 1.079124093036629
$ws.Range("N25").Value = 1.069984719925969
